# Xcg/Ycg ESTIMATION METHOD COMPARISON tables: the SFORZA and TORENBEEK_1982
# method rows were reordered (swapped) in the FUSELAGE and WING sheets.
#
# FUSELAGE: Xcg comparison table -> rows 23 (SFORZA) / 24 (TORENBEEK_1982) swap
# WING:     Xcg comparison table -> rows 23 (SFORZA) / 24 (TORENBEEK_1982) swap
#           Ycg comparison table -> rows 27 (SFORZA) / 28 (TORENBEEK_1982) swap

$wb = $excel.ActiveWorkbook

function Swap-Row {
    param($ws, [int]$row1, [int]$row2)

    $labelCell1 = $ws.Cells.Item($row1, 1)
    $valueCell1 = $ws.Cells.Item($row1, 3)
    $labelCell2 = $ws.Cells.Item($row2, 1)
    $valueCell2 = $ws.Cells.Item($row2, 3)

    $label1 = $labelCell1.Value2
    $value1 = $valueCell1.Value2
    $label2 = $labelCell2.Value2
    $value2 = $valueCell2.Value2

    $labelCell1.Value = $label2
    $valueCell1.Value = $value2
    $labelCell2.Value = $label1
    $valueCell2.Value = $value1
}

$wsFuselage = $wb.Worksheets.Item("FUSELAGE")
Swap-Row $wsFuselage 23 24

$wsWing = $wb.Worksheets.Item("WING")
Swap-Row $wsWing 23 24
Swap-Row $wsWing 27 28
